$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item("Google Shape;547;p42")

# Resize the textbox (offset unchanged; height 415500 -> 620652)
$shp.Height = 620652

# Update the run text: en dash + right single quotation mark characters
$tr = $shp.TextFrame.TextRange
$tr.Text = "0 – Doesn’t have a name    1 - Has a name"
